# Update "想去人数" (F) and "最低票价" (G) figures that changed between
# the previous and newly generated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value  = 1789
$ws1.Range("F4").Value  = 455
$ws1.Range("F8").Value  = 337
$ws1.Range("F9").Value  = 1737
$ws1.Range("F10").Value = 367
$ws1.Range("F13").Value = 338
$ws1.Range("F14").Value = 682
$ws1.Range("F15").Value = 12811
$ws1.Range("G15").Value = 70.2
$ws1.Range("F16").Value = 12808
$ws1.Range("F17").Value = 958
$ws1.Range("F21").Value = 52
$ws1.Range("F22").Value = 568
$ws1.Range("F23").Value = 2006
$ws1.Range("F27").Value = 37
$ws1.Range("F28").Value = 251

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F5").Value  = 1789
$ws4.Range("F6").Value  = 455
$ws4.Range("F13").Value = 337
$ws4.Range("F14").Value = 1737
$ws4.Range("F15").Value = 367
$ws4.Range("F18").Value = 338
$ws4.Range("F20").Value = 682
$ws4.Range("F21").Value = 12811
$ws4.Range("G21").Value = 70.2
$ws4.Range("F22").Value = 12809
$ws4.Range("F23").Value = 958
$ws4.Range("F27").Value = 52
$ws4.Range("F28").Value = 568
$ws4.Range("F31").Value = 2006
$ws4.Range("F37").Value = 37
$ws4.Range("F38").Value = 251
